$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("M2").Value = "ليكر ان مكسور"
$ws.Range("M4").Value = ""
